# Remove "Non classée" as an account type: change row 13 (account "Non classé")
# on the "Comptes" sheet to use Type "Dépenses" instead, and update its Numéro.

$wb = $excel.ActiveWorkbook
$wsComptes = $wb.Worksheets.Item("Comptes")

$wsComptes.Range("B13").Value = "Dépenses"
$wsComptes.Range("E13").Value = 5999

# Make "Comptes" the active/selected sheet with A13 selected.
$wsComptes.Activate()
$wsComptes.Range("A13").Select()
